$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "어느 대기업 계열사에서 온 AI글 요청"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/a-conglomerate-subsidiary-request/#utm_source=rss&utm_medium=rss&utm_campaign=a-conglomerate-subsidiary-request"

$ws.Range("D23").Value = "Free Data Science Courses offered by Kaggle | Kaggle Free courses | Machine Learning | Data Magic"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2769"

$ws.Range("D35").Value = "FCM 모바일 푸쉬 성능 개선"
$ws.Range("E35").Value = "http://docs.likejazz.com/fcm-push/"

$ws.Range("D43").Value = "jupyter notebook TOC 쓸 수 있게 extension 설치 명령어"
$ws.Range("E43").Value = "https://nittaku.tistory.com/508"

$ws.Range("D51").Value = "[python+openpose] openpose 라이브러리를 사용해서 관절 포인트 검출하기 (window 10 환경)"
$ws.Range("E51").Value = "https://bskyvision.com/1164"
